# Price tracker update: insert a new snapshot column at column B,
# shifting all existing date columns one position to the right, and
# record the new snapshot's timestamp header. The new column's price
# cells are left blank (no scrape recorded yet for this run), matching
# the source data's existing convention for "no data" cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts columns B:BV to C:BW
# and carries the header row's style (s="1") onto the new column,
# matching Excel's native "insert column" behaviour.
$ws.Range("B1").EntireColumn.Insert()

# The newly inserted column has no explicit width yet; give it the same
# custom width (raw OOXML width = 21) as every other data column.
# Excel's ColumnWidth property differs from the stored raw width by the
# fixed default-font padding offset (5/6 = 0.8333333333333334), so we
# back that out to land exactly on width=21.
$ws.Columns.Item(2).ColumnWidth = 20.166666666666668

# Stamp the new snapshot's timestamp in the header row.
$ws.Range("B1").Value = "2025-12-27 00:25"
